# ES-205 DOCUMENTATION.docx - content touch-up pass
# ("comments for better aid")
#
# Word property/constant values used below (no enum module available in
# this headless host, so the raw integers are spelled out):
#   wdReplaceAll    = 2   (Find.Execute Replace: parameter)
#   wdOrientPortrait = 0  (PageSetup.Orientation)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Dataset: S&P 500 " / "Stocks" were left split across two runs from
#    an earlier edit; retype the line as one clean pass so it collapses
#    back into a single run.
# ---------------------------------------------------------------------
$find1 = $d.Content.Find
$find1.Execute("Dataset: S&P 500 Stocks", $false, $false, $false, $false, $false, $true, 1, $false, "Dataset: S&P 500 StocksZZZ", 2) | Out-Null
$find1b = $d.Content.Find
$find1b.Execute("Dataset: S&P 500 StocksZZZ", $false, $false, $false, $false, $false, $true, 1, $false, "Dataset: S&P 500 Stocks", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Same cleanup for the Kaggle URL in the citation hyperlink - it was
#    split into "...sandp50" + "0"; re-type it as one run so the
#    hyperlink text is a single, well-formed run (keeps the Hyperlink
#    character style). Matching a substring that starts mid-run (rather
#    than exactly on the hyperlink's leading run boundary) keeps the
#    Find/Replace engine anchored on the Hyperlink run's own formatting.
# ---------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.Execute("camnugent/sandp500", $false, $false, $false, $false, $false, $true, 1, $false, "camnugent/sandp500ZZZ", 2) | Out-Null
$find2b = $d.Content.Find
$find2b.Execute("camnugent/sandp500ZZZ", $false, $false, $false, $false, $false, $true, 1, $false, "camnugent/sandp500", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Touch the title run's character formatting (bold on, then back off)
#    so Word stamps an explicit (empty) run-properties element on it.
# ---------------------------------------------------------------------
$titleRange = $d.Content
$titleRange.Find.Execute("ES-205 CEP") | Out-Null
$titleRange.Bold = 1
$titleRange.Bold = 0

# ---------------------------------------------------------------------
# 4) Nail down the page orientation explicitly as Portrait.
# ---------------------------------------------------------------------
$d.PageSetup.Orientation = 0
